{"js": "// Locate the paragraph that ends the \"Similarly we could have the orders\n// being shipped...\" sentence and, right after it, insert:\n//   1. a blank paragraph\n//   2. a new paragraph describing the data engineer's responsibilities\nconst anchorText =\n  \"Similarly we could have the orders being shipped at different time periods for different customers.\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph to insert after.\");\n}\n\nconst newText =\n  \"The data engineer has to keep in mind the primary keys and the foreign keys while creating the pipeline and the data types of the keys so that the joins can be done without any issues. The data engineer needs to create 3 data sets and make sure the columns are as specified in the data model. For combining all the 3 tables together, first the engineer can join the customers and orders table and then with this combined table the shipping table can be joined. However since shipping data is available in json it needs to be first converted (I have done in Pyspark). This will give one data set where we have all the information together.\";\n\nconst blankParagraph = target.insertParagraph(\"\", Word.InsertLocation.after);\nblankParagraph.insertParagraph(newText, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Locate the paragraph that ends the \"Similarly we could have the orders\n# being shipped...\" sentence and, right after it, insert:\n#   1. a blank paragraph\n#   2. a new paragraph describing the data engineer's responsibilities\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.Execute(\"Similarly we could have the orders being shipped at different time periods for different customers.\")\n\n# Move to the end of the found text, then add a blank paragraph.\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n# Step into the newly created blank paragraph and add a second paragraph\n# after it, carrying the new sentence.\n$rng.Collapse(0)\n$rng.Move(1, 1)\n$rng.InsertParagraphAfter()\n\n$rng.Collapse(0)\n$rng.Move(1, 1)\n$rng.InsertAfter(\"The data engineer has to keep in mind the primary keys and the foreign keys while creating the pipeline and the data types of the keys so that the joins can be done without any issues. The data engineer needs to create 3 data sets and make sure the columns are as specified in the data model. For combining all the 3 tables together, first the engineer can join the customers and orders table and then with this combined table the shipping table can be joined. However since shipping data is available in json it needs to be first converted (I have done in Pyspark). This will give one data set where we have all the information together.\")\n"}
